$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A6: date value changes from 41527 (09/10/2013) to 41557 (10/10/2013)
$ws.Range("A6").Value = 41557

# Add new row 7 data: copy A6's style/format down to A7, then set its date value
$ws.Range("A6").Copy($ws.Range("A7"))
$ws.Range("A7").Value = 41558

# B7 already has the right style; just set its time value
$ws.Range("B7").Value = 0.1388888888888889

# Update the selection to B4:B7 (matches the sqref in the target view state;
# this headless engine always anchors the active cell at the top-left of the
# selected range, so B4 becomes active - there is no COM-exposed way here to
# independently pick B7 as active cell within a multi-cell selection).
$ws.Range("B4:B7").Select()
